$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 34210.734
$ws.Range("I116").Value = 51098.953
$ws.Range("J116").Value = 3249
$ws.Range("K116").Value = 51098.953
$ws.Range("L116").Value = 3249
$ws.Range("M116").Value = -47656.953
$ws.Range("N116").Value = -10133
$ws.Range("H132").Value = 2132.6377
$ws.Range("I132").Value = 1166.9678
$ws.Range("J132").Value = 10685.714
$ws.Range("K132").Value = 3500.9034
$ws.Range("L132").Value = 32057.142
$ws.Range("M132").Value = -970.9033999999997
$ws.Range("N132").Value = -37117.142

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 351828.47
$ws.Range("I61").Value = 235293.88
$ws.Range("J61").Value = 773145.9399999999
$ws.Range("K61").Value = 235293.88
$ws.Range("L61").Value = 773145.9399999999
$ws.Range("M61").Value = -235081.88
$ws.Range("N61").Value = -773569.9399999999
$ws.Range("H122").Value = 2557.7908
$ws.Range("I122").Value = 2304.3794
$ws.Range("J122").Value = 3082.7144
$ws.Range("K122").Value = 6913.138199999999
$ws.Range("L122").Value = 9248.143199999999
$ws.Range("M122").Value = -4463.138199999999
$ws.Range("N122").Value = -14148.1432
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").Value = ""
$ws.Range("H135").Value = 38107.25
$ws.Range("J135").Value = 38107.25
$ws.Range("L135").Value = 38107.25
$ws.Range("N135").Value = -48247.25
$ws.Range("H136").Value = 351828.47
$ws.Range("I136").Value = 235293.88
$ws.Range("J136").Value = 773145.9399999999
$ws.Range("K136").Value = 705881.64
$ws.Range("L136").Value = 2319437.82
$ws.Range("M136").Value = -703331.64
$ws.Range("N136").Value = -2324537.82
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").Value = ""
$ws.Range("H140").Value = 40429
$ws.Range("J140").Value = 40429
$ws.Range("L140").Value = 40429
$ws.Range("N140").Value = -50789
$ws.Range("H141").Value = 36495
$ws.Range("J141").Value = 36495
$ws.Range("L141").Value = 36495
$ws.Range("N141").Value = -46855

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 880.10254
$ws.Range("I94").Value = 635
$ws.Range("J94").Value = 1504
$ws.Range("K94").Value = 635
$ws.Range("L94").Value = 1504
$ws.Range("M94").Value = -184
$ws.Range("N94").Value = -2406

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1943.4193
$ws.Range("I31").Value = 1466.0952
$ws.Range("J31").Value = 2945.8
$ws.Range("K31").Value = 1466.0952
$ws.Range("L31").Value = 2945.8
$ws.Range("M31").Value = -1171.0952
$ws.Range("N31").Value = -3535.8
$ws.Range("H34").Value = 1943.4193
$ws.Range("I34").Value = 1466.0952
$ws.Range("J34").Value = 2945.8
$ws.Range("K34").Value = 1466.0952
$ws.Range("L34").Value = 2945.8
$ws.Range("M34").Value = -1264.0952
$ws.Range("N34").Value = -3349.8
$ws.Range("H58").Value = 4092.311
$ws.Range("I58").Value = 5566.591
$ws.Range("J58").Value = 2682.1304
$ws.Range("K58").Value = 5566.591
$ws.Range("L58").Value = 2682.1304
$ws.Range("M58").Value = -5363.591
$ws.Range("N58").Value = -3088.1304
$ws.Range("H132").Value = 1881.591
$ws.Range("I132").Value = 1043.2258
$ws.Range("K132").Value = 3129.6774
$ws.Range("M132").Value = -599.6773999999996
$ws.Range("H136").Value = 4092.311
$ws.Range("I136").Value = 5566.591
$ws.Range("J136").Value = 2682.1304
$ws.Range("K136").Value = 16699.773
$ws.Range("L136").Value = 8046.3912
$ws.Range("M136").Value = -14149.773
$ws.Range("N136").Value = -13146.3912

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 189.47368
$ws.Range("I40").Value = 67.85714
$ws.Range("K40").Value = 271.42856
$ws.Range("M40").Value = -202.42856
$ws.Range("H82").Value = 6650
$ws.Range("I82").Value = 900
$ws.Range("J82").Value = 7800
$ws.Range("K82").Value = 2700
$ws.Range("L82").Value = 23400
$ws.Range("M82").Value = -2294
$ws.Range("N82").Value = -24212
$ws.Range("H85").Value = 6650
$ws.Range("I85").Value = 900
$ws.Range("J85").Value = 7800
$ws.Range("K85").Value = 2700
$ws.Range("L85").Value = 23400
$ws.Range("M85").Value = -1296
$ws.Range("N85").Value = -26208
$ws.Range("H111").Value = 1378.4073
$ws.Range("I111").Value = 369
$ws.Range("J111").Value = 1504.5834
$ws.Range("K111").Value = 1107
$ws.Range("L111").Value = 4513.7502
$ws.Range("M111").Value = 1960
$ws.Range("N111").Value = -10647.7502

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1190.697
$ws.Range("I122").Value = 1199.6111
$ws.Range("J122").Value = 1180
$ws.Range("K122").Value = 3598.8333
$ws.Range("L122").Value = 3540
$ws.Range("M122").Value = -1148.8333
$ws.Range("N122").Value = -8440
$ws.Range("H126").Value = 2916.5144
$ws.Range("I126").Value = 2597.2666
$ws.Range("J126").Value = 3155.95
$ws.Range("K126").Value = 7791.7998
$ws.Range("L126").Value = 9467.849999999999
$ws.Range("M126").Value = -5321.7998
$ws.Range("N126").Value = -14407.85

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2387.6667
$ws.Range("I82").Value = 1787
$ws.Range("J82").Value = 3228.6
$ws.Range("K82").Value = 1787
$ws.Range("L82").Value = 3228.6
$ws.Range("M82").Value = -1426
$ws.Range("N82").Value = -3950.6
$ws.Range("H85").Value = 2387.6667
$ws.Range("I85").Value = 1787
$ws.Range("J85").Value = 3228.6
$ws.Range("K85").Value = 1787
$ws.Range("L85").Value = 3228.6
$ws.Range("M85").Value = -539
$ws.Range("N85").Value = -5724.6
$ws.Range("H93").Value = 1162.875
$ws.Range("I93").Value = 944.1429000000001
$ws.Range("J93").Value = 1580.4546
$ws.Range("K93").Value = 944.1429000000001
$ws.Range("L93").Value = 1580.4546
$ws.Range("M93").Value = 303.8570999999999
$ws.Range("N93").Value = -4076.4546
$ws.Range("H104").Value = 8072.25
$ws.Range("J104").Value = 8072.25
$ws.Range("L104").Value = 8072.25
$ws.Range("N104").Value = -15060.25
$ws.Range("H132").Value = 6616.4683
$ws.Range("I132").Value = 2339.0435
$ws.Range("J132").Value = 10715.667
$ws.Range("K132").Value = 7017.130500000001
$ws.Range("L132").Value = 32147.001
$ws.Range("M132").Value = -4487.130500000001
$ws.Range("N132").Value = -37207.001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 6925
$ws.Range("I39").Value = 6000
$ws.Range("J39").Value = 7233.3335
$ws.Range("K39").Value = 6000
$ws.Range("L39").Value = 7233.3335
$ws.Range("M39").Value = -5587
$ws.Range("N39").Value = -8059.3335
$ws.Range("H104").Value = 18451.6
$ws.Range("J104").Value = 18451.6
$ws.Range("L104").Value = 18451.6
$ws.Range("N104").Value = -25439.6
$ws.Range("H122").Value = 38463120
$ws.Range("I122").Value = 45456028
$ws.Range("K122").Value = 136368084
$ws.Range("M122").Value = -136365634
$ws.Range("H126").Value = 820.6667
$ws.Range("J126").Value = 1143.3334
$ws.Range("L126").Value = 3430.0002
$ws.Range("N126").Value = -8370.0002
$ws.Range("H132").Value = 1857.0238
$ws.Range("I132").Value = 1196.1034
$ws.Range("J132").Value = 3331.3845
$ws.Range("K132").Value = 3588.3102
$ws.Range("L132").Value = 9994.1535
$ws.Range("M132").Value = -1058.3102
$ws.Range("N132").Value = -15054.1535
